$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5514.815
$ws.Range("I40").Value = 1050
$ws.Range("J40").Value = 6790.476
$ws.Range("K40").Value = 1050
$ws.Range("L40").Value = 6790.476
$ws.Range("M40").Value = -875
$ws.Range("N40").Value = -7140.476
$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 5000
$ws.Range("K76").Value = 5000
$ws.Range("M76").Value = -4685
$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 5000
$ws.Range("K79").Value = 5000
$ws.Range("M79").Value = -3908
$ws.Range("H88").Value = 5000
$ws.Range("I88").Value = 5000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 5000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -4594
$ws.Range("N88").Value = ""
$ws.Range("H91").Value = 5000
$ws.Range("I91").Value = 5000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 5000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -3596
$ws.Range("N91").Value = ""
$ws.Range("H98").Value = 4130
$ws.Range("I98").Value = 4130
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 4130
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -2632
$ws.Range("N98").Value = ""
$ws.Range("H122").Value = 4130
$ws.Range("I122").Value = 4130
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 12390
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9940
$ws.Range("N122").Value = ""
$ws.Range("H125").Value = 1083
$ws.Range("I125").Value = 821.6667
$ws.Range("J125").Value = 1475
$ws.Range("K125").Value = 7395.0003
$ws.Range("L125").Value = 13275
$ws.Range("M125").Value = -4935.0003
$ws.Range("N125").Value = -18195
$ws.Range("H127").Value = 4000
$ws.Range("J127").Value = 4000
$ws.Range("L127").Value = 12000
$ws.Range("N127").Value = -21920
$ws.Range("H137").Value = 1618.6154
$ws.Range("I137").Value = 1459.3334
$ws.Range("K137").Value = 4378.0002
$ws.Range("M137").Value = -1828.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2333.3333
$ws.Range("I63").Value = 1000
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 1000
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -314
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 2333.3333
$ws.Range("I66").Value = 1000
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 5000
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -1568
$ws.Range("N66").Value = -21864
$ws.Range("H88").Value = 305.5
$ws.Range("I88").Value = 315
$ws.Range("J88").Value = 296
$ws.Range("K88").Value = 315
$ws.Range("L88").Value = 296
$ws.Range("M88").Value = 91
$ws.Range("N88").Value = -1108
$ws.Range("H91").Value = 305.5
$ws.Range("I91").Value = 315
$ws.Range("J91").Value = 296
$ws.Range("K91").Value = 315
$ws.Range("L91").Value = 296
$ws.Range("M91").Value = 1089
$ws.Range("N91").Value = -3104
$ws.Range("H122").Value = 9996
$ws.Range("I122").Value = 9996
$ws.Range("K122").Value = 29988
$ws.Range("M122").Value = -27538
$ws.Range("H132").Value = 1199.75
$ws.Range("I132").Value = 999
$ws.Range("K132").Value = 2997
$ws.Range("M132").Value = -467

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3037.6
$ws.Range("I20").Value = 2462.6667
$ws.Range("J20").Value = 3900
$ws.Range("K20").Value = 2462.6667
$ws.Range("L20").Value = 3900
$ws.Range("M20").Value = -2215.6667
$ws.Range("N20").Value = -4394
$ws.Range("H105").Value = 1889.8
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2699
$ws.Range("I62").Value = 2699
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2699
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2075
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 2699
$ws.Range("I65").Value = 2699
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 13495
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -10375
$ws.Range("N65").Value = ""
$ws.Range("H122").Value = 500
$ws.Range("I122").Value = 500
$ws.Range("K122").Value = 1500
$ws.Range("M122").Value = 950

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 149.63637
$ws.Range("J2").Value = 140
$ws.Range("L2").Value = 840
$ws.Range("N2").Value = -1066
$ws.Range("H4").Value = 1667156.6
$ws.Range("I4").Value = 1667156.6
$ws.Range("K4").Value = 5001469.800000001
$ws.Range("M4").Value = -5001357.800000001
$ws.Range("H47").Value = 201
$ws.Range("I47").Value = 3
$ws.Range("J47").Value = 399
$ws.Range("K47").Value = 9
$ws.Range("L47").Value = 1197
$ws.Range("M47").Value = 422
$ws.Range("N47").Value = -2059

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2999
$ws.Range("I102").Value = 2999
$ws.Range("K102").Value = 2999
$ws.Range("M102").Value = -1377
$ws.Range("H122").Value = 1400
$ws.Range("I122").Value = 1200
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3600
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1150
$ws.Range("N122").Value = -9400
$ws.Range("H126").Value = 4090.8333
$ws.Range("I126").Value = 4309
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 12927
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -10457
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 17782
$ws.Range("I5").Value = 17782
$ws.Range("K5").Value = 17782
$ws.Range("M5").Value = -17669
$ws.Range("H7").Value = 4112.25
$ws.Range("I7").Value = 4112.25
$ws.Range("K7").Value = 4112.25
$ws.Range("M7").Value = -4000.25
$ws.Range("H22").Value = 3002
$ws.Range("J22").Value = 3002
$ws.Range("L22").Value = 3002
$ws.Range("N22").Value = -3592
$ws.Range("H27").Value = 3002
$ws.Range("J27").Value = 3002
$ws.Range("L27").Value = 3002
$ws.Range("N27").Value = -3216
$ws.Range("H40").Value = 7133
$ws.Range("I40").Value = 7133
$ws.Range("K40").Value = 7133
$ws.Range("M40").Value = -6997
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = ""
$ws.Range("H126").Value = 4112.25
$ws.Range("I126").Value = 4112.25
$ws.Range("K126").Value = 12336.75
$ws.Range("M126").Value = -9866.75
$ws.Range("H132").Value = 4364.1333
$ws.Range("I132").Value = 4178.4546
$ws.Range("J132").Value = 4874.75
$ws.Range("K132").Value = 12535.3638
$ws.Range("L132").Value = 14624.25
$ws.Range("M132").Value = -10005.3638
$ws.Range("N132").Value = -19684.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8125
$ws.Range("I122").Value = 7833.3335
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 23500.0005
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -21050.0005
$ws.Range("N122").Value = -31900
$ws.Range("H126").Value = 1677.2222
$ws.Range("J126").Value = 1254
$ws.Range("L126").Value = 3762
$ws.Range("N126").Value = -8702
